$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260, shifting rows 260:356 down to 261:357.
$ws.Rows.Item(260).Insert()

# Fill in the new row's data (mirrors row 260's original record, with the
# volume/price/date columns updated).
$ws.Cells.Item(260, 1).Value = 2
$ws.Cells.Item(260, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(260, 3).Value = "Coquimbo"
$ws.Cells.Item(260, 4).Value = 44924
$ws.Cells.Item(260, 5).Value = 4
$ws.Cells.Item(260, 6).Value = 100112021
$ws.Cells.Item(260, 7).Value = "Ají"
$ws.Cells.Item(260, 8).Value = "Americana (o)"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 240
$ws.Cells.Item(260, 11).Value = 12000
$ws.Cells.Item(260, 12).Value = 13000
$ws.Cells.Item(260, 13).Value = 12500
$ws.Cells.Item(260, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(260, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(260, 16).Value = 500
$ws.Cells.Item(260, 17).Value = 25
$ws.Cells.Item(260, 18).Value = "Hortaliza"
